$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Förändrad) for rows 2-8: 45183 -> 45184 (serial date, keep formatting)
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45184
}
